# resume shorten.docx - content edit
#
# 1) Summary paragraph: isolate the word "by" into its own run (the run
#    that used to hold just the trailing period "." now holds
#    " collaborating across teams and optimizing for performance and
#    user experience." instead - the visible text is unchanged, only the
#    run boundaries move).
#
# 2) SKILLS section: merge the "Frontend: React, Vue" paragraph and the
#    "Backend: .NET, NodeJS, NestJS, ExpressJS, Laravel" paragraph into a
#    single bullet, joined with " | ", which shifts Database/Tools/Systems
#    up by one paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: split "by" out of the summary sentence into its own run
# ---------------------------------------------------------------------
$summary = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Full stack web developer*") {
        $summary = $cand
        break
    }
}

if ($summary -ne $null) {
    $pStart = $summary.Range.Start
    $pText  = $summary.Range.Text

    $byRel   = $pText.IndexOf(" by ") + 1
    $byStart = $pStart + $byRel
    $byEnd   = $byStart + 2
    # $pText includes the trailing paragraph-mark character, so this is
    # the absolute offset of the end of the paragraph (incl. the mark).
    $tailEnd = $pStart + $pText.Length

    # Re-assert the tail text ("collaborating ... experience.") as one
    # run distinct from "by" - this also absorbs the old lone-"." run.
    $rTail = $d.Range($byEnd, $tailEnd)
    $rTail.Text = " collaborating across teams and optimizing for performance and user experience."

    # Force "by" to split off from the preceding "...schedule " run by
    # toggling formatting on just that word (net no-op visually).
    $rBy = $d.Range($byStart, $byEnd)
    $rBy.Bold = 1
    $rBy.Bold = 0
}

# ---------------------------------------------------------------------
# Edit 2: merge "Frontend: React, Vue" + "Backend: ..." into one bullet
# ---------------------------------------------------------------------
$frontend = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Frontend:*") {
        $frontend = $cand
        break
    }
}

if ($frontend -ne $null) {
    $textEnd = $frontend.Range.End - 1   # position right before the paragraph mark

    $rInsert = $d.Range($textEnd, $textEnd)
    $rInsert.InsertAfter(" | ")

    # Force the new " | " separator into its own run (keeps
    # "Frontend: React, Vue" as the original, untouched run).
    $rSep = $d.Range($textEnd, $textEnd + 3)
    $rSep.Bold = 1
    $rSep.Bold = 0

    # Re-locate the paragraph (its Range.End shifted right by 3 chars)
    # and delete its trailing paragraph mark so the "Backend: ..."
    # paragraph's content joins this paragraph.
    $frontend2 = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $cand = $d.Paragraphs.Item($i)
        if ($cand.Range.Text -like "Frontend:*") {
            $frontend2 = $cand
            break
        }
    }
    $markStart = $frontend2.Range.End - 1
    $markEnd   = $frontend2.Range.End
    $rMark = $d.Range($markStart, $markEnd)
    $rMark.Delete()
}
